# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" worksheets for rows 3-5 (column F).
#   F3: 932 -> 934
#   F4: 216 -> 217
#   F5: 428 -> 429

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 934
    $ws.Range("F4").Value = 217
    $ws.Range("F5").Value = 429
}
